# Retention of Escolars no contesten: insert a new row for
# "No contesta (""escolar"")" above the "Subtotal (b)" row on the Hoja1
# report, and turn on iterative calculation (delta 1E-4) for the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Enable iterative calculation with the delta used by the updated report.
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# Push row 15 ("Subtotal (b)") and everything below it down by one row,
# copying the formatting of the row above the insertion point.
$ws.Rows("15").Insert()

# Fill in the freshly inserted row 15 with the new "escolar" bucket.
$ws.Range("A15").Value = "No contesta (""escolar"")"

# Match the author's final selection on the sheet.
$ws.Range("E15").Select() | Out-Null
